# Hortaliza, Terminal Hortofrutícola Agro Chillán - Choclo
# Weekly update: insert two new price rows (238-239) and push existing
# historical rows down by two (old 238-250 -> new 240-252).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 238, shifting rows 238:250 down to 240:252.
$ws.Rows("238:239").Insert()

# New row 238: Choclero / Primera, Región del Maule
$ws.Cells.Item(238, 1).Value = 7
$ws.Cells.Item(238, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(238, 3).Value = "Ñuble"
$ws.Cells.Item(238, 4).Value = 44931
$ws.Cells.Item(238, 5).Value = 16
$ws.Cells.Item(238, 6).Value = 100112024
$ws.Cells.Item(238, 7).Value = "Choclo"
$ws.Cells.Item(238, 8).Value = "Choclero"
$ws.Cells.Item(238, 9).Value = "Primera"
$ws.Cells.Item(238, 10).Value = 20000
$ws.Cells.Item(238, 11).Value = 250
$ws.Cells.Item(238, 12).Value = 300
$ws.Cells.Item(238, 13).Value = 275
$ws.Cells.Item(238, 14).Value = "$/unidad"
$ws.Cells.Item(238, 15).Value = "Región del Maule"
$ws.Cells.Item(238, 16).Value = 275
$ws.Cells.Item(238, 17).Value = 1
$ws.Cells.Item(238, 18).Value = "Hortaliza"

# New row 239: Choclero / Segunda, Región del Maule
$ws.Cells.Item(239, 1).Value = 7
$ws.Cells.Item(239, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(239, 3).Value = "Ñuble"
$ws.Cells.Item(239, 4).Value = 44931
$ws.Cells.Item(239, 5).Value = 16
$ws.Cells.Item(239, 6).Value = 100112024
$ws.Cells.Item(239, 7).Value = "Choclo"
$ws.Cells.Item(239, 8).Value = "Choclero"
$ws.Cells.Item(239, 9).Value = "Segunda"
$ws.Cells.Item(239, 10).Value = 15000
$ws.Cells.Item(239, 11).Value = 200
$ws.Cells.Item(239, 12).Value = 200
$ws.Cells.Item(239, 13).Value = 200
$ws.Cells.Item(239, 14).Value = "$/unidad"
$ws.Cells.Item(239, 15).Value = "Región del Maule"
$ws.Cells.Item(239, 16).Value = 200
$ws.Cells.Item(239, 17).Value = 1
$ws.Cells.Item(239, 18).Value = "Hortaliza"
